$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly activity rows appended after the existing last row (207).
$newRows = @(
    @{ Row = 208; A = "Steven";   B = 45486; C = "Walk";    D = 33; E = 1.51;  F = 59; G = 33; H = 0;  I = 0;  J = 0; K = 0; L = "Brave Leopard"; M = 5 },
    @{ Row = 209; A = "Eric";     B = 45486; C = "Workout"; D = 54; E = 0;     F = 0;  G = 6;  H = 38; I = 10; J = 3; K = 0; L = "Wily Hyena";    M = 5 },
    @{ Row = 210; A = "Jeremiah"; B = 45487; C = "Ride";    D = 30; E = 11.12; F = 0;  G = 2;  H = 23; I = 5;  J = 0; K = 0; L = "Wily Hyena";    M = 5 },
    @{ Row = 211; A = "Jeremiah"; B = 45487; C = "Workout"; D = 18; E = 0;     F = 0;  G = 18; H = 0;  I = 0;  J = 0; K = 0; L = "Wily Hyena";    M = 5 },
    @{ Row = 212; A = "Steven";   B = 45487; C = "Walk";    D = 42; E = 1.99;  F = 98; G = 41; H = 1;  I = 0;  J = 0; K = 0; L = "Brave Leopard"; M = 5 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy the date cell's formatting from the row above so the new date
    # cell reuses the existing date style instead of minting a new one.
    $ws.Range("B" + ($row - 1)).Copy() | Out-Null
    $ws.Range("B" + $row).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $row).Value = $r.A
    $ws.Range("B" + $row).Value = $r.B
    $ws.Range("C" + $row).Value = $r.C
    $ws.Range("D" + $row).Value = $r.D
    $ws.Range("E" + $row).Value = $r.E
    $ws.Range("F" + $row).Value = $r.F
    $ws.Range("G" + $row).Value = $r.G
    $ws.Range("H" + $row).Value = $r.H
    $ws.Range("I" + $row).Value = $r.I
    $ws.Range("J" + $row).Value = $r.J
    $ws.Range("K" + $row).Value = $r.K
    $ws.Range("L" + $row).Value = $r.L
    $ws.Range("M" + $row).Value = $r.M
}

$ws.Range("A213").Select() | Out-Null
